# "Generate Report for Archive" — refresh the localization-status report:
#   * flip the in-flight status label from "Ready for handoff" to "In Translation"
#     on all three sheets (Overview's zh-cn/de-de columns, and the Status column
#     on each per-locale sheet)
#   * re-fit the Status column now that the shorter label no longer needs as
#     much room

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$null = $overview.Cells.Replace("Ready for handoff", "In Translation")
$null = $zhcn.Cells.Replace("Ready for handoff", "In Translation")
$null = $dede.Cells.Replace("Ready for handoff", "In Translation")

# The Status columns can now be narrower since "In Translation" is shorter
# than "Ready for handoff". Resize them to match the new content width.
$overview.Columns.Item(5).ColumnWidth = 12.5   # Overview!E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # Overview!F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # de-de!C (Status)
